$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 26.5325808198159
$ws.Range("C2").Value = 10.87194340893262
$ws.Range("E2").Value = 8.508707568532952
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.754285107724709
$ws.Range("L2").Value = 10.10588601429479
$ws.Range("M2").Value = 19.99784917743138
$ws.Range("N2").Value = 22.78493605613365

# Row 3
$ws.Range("B3").Value = 26.15954974725561
$ws.Range("C3").Value = 10.35922384877414
$ws.Range("E3").Value = 8.468465003891037
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.759250733307188
$ws.Range("L3").Value = 10.12075865247675
$ws.Range("M3").Value = 19.93572695825971
$ws.Range("N3").Value = 22.81166820434175

# Row 4
$ws.Range("B4").Value = 25.93638612204292
$ws.Range("C4").Value = 10.03496306661431
$ws.Range("E4").Value = 8.443150268554952
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.762452321319787
$ws.Range("L4").Value = 10.1314124704953
$ws.Range("M4").Value = 19.90237546946706
$ws.Range("N4").Value = 22.82988578879255

# Row 5
$ws.Range("B5").Value = 25.84703341556853
$ws.Range("C5").Value = 9.900656208878161
$ws.Range("E5").Value = 8.432682140860797
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.763795556619814
$ws.Range("L5").Value = 10.13613645928734
$ws.Range("M5").Value = 19.88999585225228
$ws.Range("N5").Value = 22.83776156933077

# Row 6
$ws.Range("B6").Value = 25.83229553801356
$ws.Range("C6").Value = 9.878230369426666
$ws.Range("E6").Value = 8.430934729211796
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.764020933892088
$ws.Range("L6").Value = 10.13694396661659
$ws.Range("M6").Value = 19.88801355277494
$ws.Range("N6").Value = 22.8390965861841

# Row 7
$ws.Range("B7").Value = 25.93517450658805
$ws.Range("C7").Value = 10.03316023909129
$ws.Range("E7").Value = 8.443009707693818
$ws.Range("F7").Value = 15.26647399323728
$ws.Range("G7").Value = 3.762470280295923
$ws.Range("L7").Value = 10.13147463157878
$ws.Range("M7").Value = 19.9022036012156
$ws.Range("N7").Value = 22.8299901764831

# Row 8
$ws.Range("B8").Value = 26.40280652250587
$ws.Range("C8").Value = 10.69723484619459
$ws.Range("E8").Value = 8.494957701935297
$ws.Range("F8").Value = 16.53996406344769
$ws.Range("G8").Value = 3.755965671300774
$ws.Range("L8").Value = 10.11069808518524
$ws.Range("M8").Value = 19.97543916197705
$ws.Range("N8").Value = 22.79377790053795

# Row 9
$ws.Range("B9").Value = 27.36150287516772
$ws.Range("C9").Value = 11.91677023926329
$ws.Range("E9").Value = 8.592019561776347
$ws.Range("F9").Value = 19.0027458068253
$ws.Range("G9").Value = 3.744413731916017
$ws.Range("L9").Value = 10.08204144435734
$ws.Range("M9").Value = 20.15673335409972
$ws.Range("N9").Value = 22.73715384704517

# Row 10
$ws.Range("B10").Value = 28.08448085363
$ws.Range("C10").Value = 12.75372766523201
$ws.Range("E10").Value = 8.6604010657814
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.736649289934628
$ws.Range("L10").Value = 10.06836832274663
$ws.Range("M10").Value = 20.31236146263987
$ws.Range("N10").Value = 22.70442739822353

# Row 11
$ws.Range("B11").Value = 28.4159816863863
$ws.Range("C11").Value = 13.12028078976924
$ws.Range("E11").Value = 8.69087387557674
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.733271629457257
$ws.Range("L11").Value = 10.063752710096
$ws.Range("M11").Value = 20.38789035102235
$ws.Range("N11").Value = 22.69148904841749

# Row 12
$ws.Range("B12").Value = 28.54176936896631
$ws.Range("C12").Value = 13.2569545653075
$ws.Range("E12").Value = 8.702321835748926
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.732014619839113
$ws.Range("L12").Value = 10.06223565608888
$ws.Range("M12").Value = 20.41715788013508
$ws.Range("N12").Value = 22.68687177828972

# Row 13
$ws.Range("B13").Value = 28.51466926711071
$ws.Range("C13").Value = 13.22761554811175
$ws.Range("E13").Value = 8.699860382663912
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.73228436193271
$ws.Range("L13").Value = 10.0625521162573
$ws.Range("M13").Value = 20.4108251946902
$ws.Range("N13").Value = 22.68785360815004

# Row 14
$ws.Range("B14").Value = 28.42632587468568
$ws.Range("C14").Value = 13.13156822490637
$ws.Range("E14").Value = 8.691817539832586
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.733167773830303
$ws.Range("L14").Value = 10.06362327576968
$ws.Range("M14").Value = 20.39028493836955
$ws.Range("N14").Value = 22.69110351600294

# Row 15
$ws.Range("B15").Value = 28.37224273963597
$ws.Range("C15").Value = 13.072456348462
$ws.Range("E15").Value = 8.686879152740797
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.733711753995558
$ws.Range("L15").Value = 10.06430944735781
$ws.Range("M15").Value = 20.37778977975461
$ws.Range("N15").Value = 22.6931309868248

# Row 16
$ws.Range("B16").Value = 28.06285875329699
$ws.Range("C16").Value = 12.7294788932259
$ws.Range("E16").Value = 8.658396813539641
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.736873120954865
$ws.Range("L16").Value = 10.06870225513717
$ws.Range("M16").Value = 20.30751952137569
$ws.Range("N16").Value = 22.70531233468944

# Row 17
$ws.Range("B17").Value = 27.87364489726824
$ws.Range("C17").Value = 12.51536983570059
$ws.Range("E17").Value = 8.640761207248111
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.738851948608838
$ws.Range("L17").Value = 10.0718080836572
$ws.Range("M17").Value = 20.26561345310355
$ws.Range("N17").Value = 22.71328576975852

# Row 18
$ws.Range("B18").Value = 27.76506726585598
$ws.Range("C18").Value = 12.39088943308716
$ws.Range("E18").Value = 8.63055798476519
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.740004663072678
$ws.Range("L18").Value = 10.07374548280177
$ws.Range("M18").Value = 20.24195659666517
$ws.Range("N18").Value = 22.71805521911681

# Row 19
$ws.Range("B19").Value = 27.72835193747437
$ws.Range("C19").Value = 12.34851706580842
$ws.Range("E19").Value = 8.627093104324107
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.740397455660881
$ws.Range("L19").Value = 10.07442738566657
$ws.Range("M19").Value = 20.23402387630836
$ws.Range("N19").Value = 22.71970150299151

# Row 20
$ws.Range("B20").Value = 27.89376171310379
$ws.Range("C20").Value = 12.53830053549523
$ws.Range("E20").Value = 8.642644726767092
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 3.738639794883552
$ws.Range("L20").Value = 10.0714618333898
$ws.Range("M20").Value = 20.27002832364998
$ws.Range("N20").Value = 22.71241799503484

# Row 21
$ws.Range("B21").Value = 28.4522684831659
$ws.Range("C21").Value = 13.15983817169666
$ws.Range("E21").Value = 8.69418240092053
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 3.732907697606802
$ws.Range("L21").Value = 10.06330238666582
$ws.Range("M21").Value = 20.39630014295758
$ws.Range("N21").Value = 22.69014126497528

# Row 22
$ws.Range("B22").Value = 28.81872037302634
$ws.Range("C22").Value = 13.55359006523871
$ws.Range("E22").Value = 8.727332217326721
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 3.729289811334968
$ws.Range("L22").Value = 10.05931486709932
$ws.Range("M22").Value = 20.48270288054509
$ws.Range("N22").Value = 22.67722792663474

# Row 23
$ws.Range("B23").Value = 28.62304585207984
$ws.Range("C23").Value = 13.34460386689089
$ws.Range("E23").Value = 8.709688381498351
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 3.731209055193443
$ws.Range("L23").Value = 10.0613199889217
$ws.Range("M23").Value = 20.43623836212789
$ws.Range("N23").Value = 22.68396878830435

# Row 24
$ws.Range("B24").Value = 27.88466625628743
$ws.Range("C24").Value = 12.52793787681247
$ws.Range("E24").Value = 8.641793387785592
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 3.738735662661301
$ws.Range("L24").Value = 10.0716179002659
$ws.Range("M24").Value = 20.26803100312557
$ws.Range("N24").Value = 22.71280973860588

# Row 25
$ws.Range("B25").Value = 27.09838626738275
$ws.Range("C25").Value = 11.59668655509079
$ws.Range("E25").Value = 8.56627633706456
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 3.74741111976862
$ws.Range("L25").Value = 10.08849840799328
$ws.Range("M25").Value = 20.10370514512346
$ws.Range("N25").Value = 22.75092111551031
